# ADD EVENT ON EXCEL
# Adds a new "IsOnBase" column to the Event sheet, nine new event rows,
# and a new "HAS_CHILDREN" condition row on the EventCondition sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Event
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Event")

# New column M: IsOnBase
$ws.Range("M1").Value = "IsOnBase"
$ws.Range("M2").Value = "bool"
$ws.Range("M3").Value = "是否在初始事件库"
$ws.Range("M4").Value = $true
$ws.Range("M5").Value = $true

# New event rows 6-14
$ws.Range("B6").Value = "佳偶天成"
$ws.Range("C6").Value = "你找到了配偶"
$ws.Range("D6").Value = "Normal"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = $false
$ws.Range("I6").Value = $false
$ws.Range("L6").Value = "MARRY,0,0"
$ws.Range("M6").Value = $true

$ws.Range("B7").Value = "意外身故"
$ws.Range("C7").Value = "你找到了配偶"
$ws.Range("D7").Value = "Normal"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = $false
$ws.Range("L7").Value = "DIE,0,0"
$ws.Range("M7").Value = $true

$ws.Range("B8").Value = "痛失爱子"
$ws.Range("C8").Value = "你找到了配偶"
$ws.Range("D8").Value = "Normal"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = $false
$ws.Range("L8").Value = "DIE_SON,0,0"
$ws.Range("M8").Value = $true

$ws.Range("B9").Value = "明镜高悬"
$ws.Range("C9").Value = "你找到了配偶"
$ws.Range("D9").Value = "Normal"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = $false
$ws.Range("I9").Value = $false
$ws.Range("L9").Value = "ADD_RESOURCE,0,0"
$ws.Range("M9").Value = $true

$ws.Range("B10").Value = "偶得机缘"
$ws.Range("C10").Value = "你找到了配偶"
$ws.Range("D10").Value = "Normal"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = $false
$ws.Range("I10").Value = $false
$ws.Range("L10").Value = "ADD_CURR_EVENT,11,1"
$ws.Range("M10").Value = $true

$ws.Range("B11").Value = "初窥门径"
$ws.Range("C11").Value = "你找到了配偶"
$ws.Range("D11").Value = "Normal"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = $false
$ws.Range("I11").Value = $false
$ws.Range("L11").Value = "ADD_CURR_EVENT,14,1"
$ws.Range("M11").Value = $false

$ws.Range("B12").Value = "登堂入室"
$ws.Range("C12").Value = "你找到了配偶"
$ws.Range("D12").Value = "Normal"
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = $false
$ws.Range("I12").Value = $false
$ws.Range("L12").Value = "ADD_MONEY,100,0"
$ws.Range("M12").Value = $true

$ws.Range("B13").Value = "小有进财"
$ws.Range("C13").Value = "你找到了配偶"
$ws.Range("D13").Value = "Normal"
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("L13").Value = "ADD_MONEY,10,0"
$ws.Range("M13").Value = $true

$ws.Range("B14").Value = "富甲一方"
$ws.Range("C14").Value = "你找到了配偶"
$ws.Range("D14").Value = "Normal"
$ws.Range("E14").Value = 11
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = $false
$ws.Range("I14").Value = $false
$ws.Range("L14").Value = "ADD_MONEY,100,0"
$ws.Range("M14").Value = $false

$ws.Range("J14").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: EventCondition
# ---------------------------------------------------------------
$wsCond = $wb.Worksheets.Item("EventCondition")

$wsCond.Range("B4").Value = "AGE"

$wsCond.Range("B5").Value = "HAS_CHILDREN"
$wsCond.Range("C5").Value = 0
$wsCond.Range("D5").Value = 0

$wsCond.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet: EventEffect (selection only)
# ---------------------------------------------------------------
$wsEffect = $wb.Worksheets.Item("EventEffect")
$wsEffect.Range("B4").Select() | Out-Null

# Return focus to the Event sheet, which was the active tab originally.
$ws.Activate() | Out-Null
$ws.Range("J14").Select() | Out-Null
